$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("老師上課 + 監考時數")

# Update D column (上課時數) values for rows where it changed from 0
$ws.Range("D2").Value = 455
$ws.Range("D7").Value = 140
$ws.Range("D11").Value = 105
$ws.Range("D14").Value = 35
$ws.Range("D20").Value = 350
$ws.Range("D27").Value = 35
$ws.Range("D29").Value = 350
$ws.Range("D33").Value = 140
$ws.Range("D37").Value = 140
$ws.Range("D40").Value = 140
$ws.Range("D44").Value = 70
$ws.Range("D46").Value = 385
$ws.Range("D50").Value = 350
$ws.Range("D54").Value = 210
$ws.Range("D57").Value = 35
$ws.Range("D59").Value = 210
$ws.Range("D64").Value = 140
$ws.Range("D68").Value = 105
$ws.Range("D70").Value = 385
$ws.Range("D79").Value = 105
$ws.Range("D82").Value = 105
$ws.Range("D85").Value = 350
$ws.Range("D90").Value = 350
$ws.Range("D96").Value = 245
$ws.Range("D99").Value = 210
$ws.Range("D102").Value = 140
$ws.Range("D105").Value = 105
$ws.Range("D118").Value = 35
$ws.Range("D125").Value = 280
$ws.Range("D133").Value = 70
$ws.Range("D135").Value = 175
$ws.Range("D138").Value = 175
$ws.Range("D142").Value = 105
$ws.Range("D145").Value = 280
$ws.Range("D151").Value = 70
$ws.Range("D154").Value = 315
$ws.Range("D158").Value = 140

# Update G column (平均) with new formula for existing data rows 2-167
$formulaRows = @(2,7,11,14,16,18,20,24,26,27,29,33,36,37,40,43,44,46,50,54,57,59,64,68,70,76,77,79,82,85,90,94,96,99,102,105,108,110,112,114,116,118,120,121,122,123,125,130,132,133,135,138,142,145,149,151,154,158,161,162,163,164,166,167)
foreach ($r in $formulaRows) {
    $ws.Range("G$r").Formula = "=((SUM(`$D`$2:`$D`$167*`$C`$2:`$C`$167)-F161-F162-F163-F164+SUM(`$E`$2:`$E`$167))/SUM(`$C`$2:`$C`$167))*C$r"
}

# Add new row 168 with G168 formula
$ws.Range("G168").Formula = "=((SUM(`$D`$2:`$D`$167*`$C`$2:`$C`$167)-F161-F162-F163-F164+SUM(`$E`$2:`$E`$167))/SUM(`$C`$2:`$C`$167))*C168"
